$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the "batsman" column (D,E) for ownTeam/oppTeam,
# shifting the old batsman..sr columns (D:I) to F:K.
$ws.Columns("D:E").Insert()

# --- Header row ---
$ws.Cells.Item(1, 1).Value = 'venue'
$ws.Cells.Item(1, 2).Value = 'date'
$ws.Cells.Item(1, 3).Value = 'result'
$ws.Cells.Item(1, 4).Value = 'ownTeam'
$ws.Cells.Item(1, 5).Value = 'oppTeam'
$ws.Cells.Item(1, 6).Value = 'batsman'
$ws.Cells.Item(1, 7).Value = 'totalRuns'
$ws.Cells.Item(1, 8).Value = 'totalBalls'
$ws.Cells.Item(1, 9).Value = 'total4s'
$ws.Cells.Item(1, 10).Value = 'total6s'
$ws.Cells.Item(1, 11).Value = 'sr'

# --- Data rows ---
$dataRows = @(
    @(' Abu Dhabi', ' October 25 2020', 'Royals won by 8 wickets (with 10 balls remaining)', 'Mumbai Indians', 'Rajasthan Royals', 'Ishan Kishan ', '37', '36', '4', '1', '102.77'),
    @(' Abu Dhabi', ' October 28 2020', 'Mumbai won by 5 wickets (with 5 balls remaining)', 'Mumbai Indians', 'Royal Challengers Bangalore', 'Ishan Kishan ', '25', '19', '3', '1', '131.57'),
    @(' Sharjah', ' November 03 2020', 'Sunrisers won by 10 wickets (with 17 balls remaining)', 'Mumbai Indians', 'Sunrisers Hyderabad', 'Ishan Kishan ', '33', '30', '1', '2', '110.00'),
    @(' Dubai (DSC)', ' November 05 2020', 'Mumbai won by 57 runs', 'Mumbai Indians', 'Delhi Capitals', 'Ishan Kishan ', '55', '30', '4', '3', '183.33'),
    @(' Dubai (DSC)', ' November 10 2020', 'Mumbai won by 5 wickets (with 8 balls remaining)', 'Mumbai Indians', 'Delhi Capitals', 'Ishan Kishan ', '33', '19', '3', '1', '173.68'),
    @(' Sharjah', ' October 23 2020', 'Mumbai won by 10 wickets (with 46 balls remaining)', 'Mumbai Indians', 'Chennai Super Kings', 'Ishan Kishan ', '68', '37', '6', '5', '183.78'),
    @(' Dubai (DSC)', ' October 31 2020', 'Mumbai won by 9 wickets (with 34 balls remaining)', 'Mumbai Indians', 'Delhi Capitals', 'Ishan Kishan ', '72', '47', '8', '3', '153.19'),
    @(' Abu Dhabi', ' October 01 2020', 'Mumbai won by 48 runs', 'Mumbai Indians', 'Kings XI Punjab', 'Ishan Kishan ', '28', '32', '1', '1', '87.50'),
    @(' Abu Dhabi', ' October 06 2020', 'Mumbai won by 57 runs', 'Mumbai Indians', 'Rajasthan Royals', 'Ishan Kishan ', '0', '1', '0', '0', '0.00'),
    @(' Dubai (DSC)', ' October 18 2020', 'Match tied (Kings XI won the one-over eliminator)', 'Mumbai Indians', 'Kings XI Punjab', 'Ishan Kishan ', '7', '7', '1', '0', '100.00'),
    @(' Dubai (DSC)', ' September 28 2020', 'Match tied (RCB won the one-over eliminator)', 'Mumbai Indians', 'Royal Challengers Bangalore', 'Ishan Kishan ', '99', '58', '2', '9', '170.68'),
    @(' Abu Dhabi', ' October 11 2020', 'Mumbai won by 5 wickets (with 2 balls remaining)', 'Mumbai Indians', 'Delhi Capitals', 'Ishan Kishan ', '28', '15', '2', '2', '186.66'),
    @(' Sharjah', ' October 04 2020', 'Mumbai won by 34 runs', 'Mumbai Indians', 'Sunrisers Hyderabad', 'Ishan Kishan ', '31', '23', '1', '2', '134.78')
)

# Columns G:K (totalRuns,totalBalls,total4s,total6s,sr) hold numbers recorded as
# text in the source data (matches the sheet's numberStoredAsText markers), so
# force a text quote-prefix on entry for those columns to avoid silent numeric coercion.
$textCols = @(7, 8, 9, 10, 11)

$r = 2
foreach ($row in $dataRows) {
    $c = 1
    foreach ($val in $row) {
        if ($textCols -contains $c) {
            $ws.Cells.Item($r, $c).Value = "'" + $val
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
        $c++
    }
    $r++
}

